$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 720.4386
$ws.Range("J17").Value = 720.4386
$ws.Range("L17").Value = 2161.3158
$ws.Range("N17").Value = -2497.3158

$ws.Range("H28").Value = 1452.8422
$ws.Range("J28").Value = 1676.5555
$ws.Range("L28").Value = 1676.5555
$ws.Range("N28").Value = -2646.5555

$ws.Range("H82").Value = 9862.666999999999
$ws.Range("J82").Value = 9925
$ws.Range("L82").Value = 29775
$ws.Range("N82").Value = -30587

$ws.Range("H85").Value = 9862.666999999999
$ws.Range("J85").Value = 9925
$ws.Range("L85").Value = 29775
$ws.Range("N85").Value = -32583

$ws.Range("H93").Value = 52500
$ws.Range("J93").Value = 52500
$ws.Range("L93").Value = 52500
$ws.Range("N93").Value = -57492

$ws.Range("H132").Value = 99746.05
$ws.Range("I132").Value = 221865.47
$ws.Range("K132").Value = 665596.41
$ws.Range("M132").Value = -663066.41

$ws.Range("H137").Value = 1870.8182
$ws.Range("I137").Value = 1529.0454
$ws.Range("J137").Value = 2554.3635
$ws.Range("K137").Value = 4587.1362
$ws.Range("L137").Value = 7663.0905
$ws.Range("M137").Value = -2037.1362
$ws.Range("N137").Value = -12763.0905

$ws.Range("H138").Value = 5830.9443
$ws.Range("I138").Value = 907.4
$ws.Range("J138").Value = 6949.9316
$ws.Range("K138").Value = 2722.2
$ws.Range("L138").Value = 20849.7948
$ws.Range("M138").Value = 2417.8
$ws.Range("N138").Value = -31129.7948

$ws.Range("H140").Value = 60838.89
$ws.Range("J140").Value = 59693.75
$ws.Range("L140").Value = 59693.75
$ws.Range("N140").Value = -70053.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3962.524
$ws.Range("I32").Value = 2225.0981
$ws.Range("K32").Value = 2225.0981
$ws.Range("M32").Value = -1938.0981

$ws.Range("H61").Value = 25713.766
$ws.Range("I61").Value = 28646.143
$ws.Range("J61").Value = 23661.1
$ws.Range("K61").Value = 28646.143
$ws.Range("L61").Value = 23661.1
$ws.Range("M61").Value = -28434.143
$ws.Range("N61").Value = -24085.1

$ws.Range("H74").Value = 5322627
$ws.Range("I74").Value = 8929681
$ws.Range("J74").Value = 6967.7896
$ws.Range("K74").Value = 8929681
$ws.Range("L74").Value = 6967.7896
$ws.Range("M74").Value = -8928807
$ws.Range("N74").Value = -8715.7896

$ws.Range("H77").Value = 5322627
$ws.Range("I77").Value = 8929681
$ws.Range("J77").Value = 6967.7896
$ws.Range("K77").Value = 44648405
$ws.Range("L77").Value = 34838.948
$ws.Range("M77").Value = -44644037
$ws.Range("N77").Value = -43574.948

$ws.Range("H102").Value = 763023.5
$ws.Range("I102").Value = 858025.75
$ws.Range("K102").Value = 858025.75
$ws.Range("M102").Value = -856403.75

$ws.Range("H132").Value = 29150.621
$ws.Range("I132").Value = 31018.45
$ws.Range("K132").Value = 93055.35000000001
$ws.Range("M132").Value = -90525.35000000001

$ws.Range("H136").Value = 25713.766
$ws.Range("I136").Value = 28646.143
$ws.Range("J136").Value = 23661.1
$ws.Range("K136").Value = 85938.429
$ws.Range("L136").Value = 70983.29999999999
$ws.Range("M136").Value = -83388.429
$ws.Range("N136").Value = -76083.29999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 111115580
$ws.Range("I105").Value = 142862320
$ws.Range("K105").Value = 142862320
$ws.Range("M105").Value = -142860573

$ws.Range("H107").Value = 17097468
$ws.Range("I107").Value = 20205428
$ws.Range("K107").Value = 20205428
$ws.Range("M107").Value = -20203508

$ws.Range("H132").Value = 81297.44500000001
$ws.Range("I132").Value = 54410
$ws.Range("K132").Value = 54410
$ws.Range("M132").Value = -49350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 18871920
$ws.Range("I31").Value = 43479944
$ws.Range("J31").Value = 5768.3335
$ws.Range("K31").Value = 43479944
$ws.Range("L31").Value = 5768.3335
$ws.Range("M31").Value = -43479649
$ws.Range("N31").Value = -6358.3335

$ws.Range("H34").Value = 18871920
$ws.Range("I34").Value = 43479944
$ws.Range("J34").Value = 5768.3335
$ws.Range("K34").Value = 43479944
$ws.Range("L34").Value = 5768.3335
$ws.Range("M34").Value = -43479742
$ws.Range("N34").Value = -6172.3335

$ws.Range("H58").Value = 3335315.8
$ws.Range("J58").Value = 3200
$ws.Range("L58").Value = 3200
$ws.Range("N58").Value = -3606

$ws.Range("H94").Value = 2171.4
$ws.Range("I94").Value = 337.5
$ws.Range("K94").Value = 337.5
$ws.Range("M94").Value = 113.5

$ws.Range("H136").Value = 3335315.8
$ws.Range("J136").Value = 3200
$ws.Range("L136").Value = 9600
$ws.Range("N136").Value = -14700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 124.478264
$ws.Range("J2").Value = 139.625
$ws.Range("L2").Value = 837.75
$ws.Range("N2").Value = -1063.75

$ws.Range("H5").Value = 738.2308
$ws.Range("I5").Value = 461.625
$ws.Range("J5").Value = 1180.8
$ws.Range("K5").Value = 1384.875
$ws.Range("L5").Value = 3542.4
$ws.Range("M5").Value = -1272.875
$ws.Range("N5").Value = -3766.4

$ws.Range("H7").Value = 212.71428
$ws.Range("I7").Value = 219.8
$ws.Range("J7").Value = 195
$ws.Range("K7").Value = 659.4000000000001
$ws.Range("L7").Value = 585
$ws.Range("M7").Value = -547.4000000000001
$ws.Range("N7").Value = -809

$ws.Range("H22").Value = 827
$ws.Range("I22").Value = 633.75
$ws.Range("K22").Value = 1901.25
$ws.Range("M22").Value = -1732.25

$ws.Range("H27").Value = 827
$ws.Range("I27").Value = 633.75
$ws.Range("K27").Value = 1901.25
$ws.Range("M27").Value = -1799.25

$ws.Range("H68").Value = 155807.22
$ws.Range("I68").Value = 1999.75
$ws.Range("J68").Value = 168897.22
$ws.Range("K68").Value = 5999.25
$ws.Range("L68").Value = 506691.66
$ws.Range("M68").Value = -5188.25
$ws.Range("N68").Value = -508313.66

$ws.Range("H71").Value = 155807.22
$ws.Range("I71").Value = 1999.75
$ws.Range("J71").Value = 168897.22
$ws.Range("K71").Value = 17997.75
$ws.Range("L71").Value = 1520074.98
$ws.Range("M71").Value = -13941.75
$ws.Range("N71").Value = -1528186.98

$ws.Range("H135").Value = 738.2308
$ws.Range("I135").Value = 461.625
$ws.Range("J135").Value = 1180.8
$ws.Range("K135").Value = 4154.625
$ws.Range("L135").Value = 10627.2
$ws.Range("M135").Value = -1619.625
$ws.Range("N135").Value = -15697.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 800
$ws.Range("I29").Value = 800
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 800
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -510
$ws.Range("N29").ClearContents()

$ws.Range("H70").Value = 5960692
$ws.Range("I70").Value = 23813524
$ws.Range("K70").Value = 23813524
$ws.Range("M70").Value = -23813254

$ws.Range("H73").Value = 5960692
$ws.Range("I73").Value = 23813524
$ws.Range("K73").Value = 23813524
$ws.Range("M73").Value = -23812588

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2940.923
$ws.Range("I22").Value = 2229.4167
$ws.Range("J22").Value = 3550.7856
$ws.Range("K22").Value = 2229.4167
$ws.Range("L22").Value = 3550.7856
$ws.Range("M22").Value = -1934.4167
$ws.Range("N22").Value = -4140.7856

$ws.Range("H27").Value = 2940.923
$ws.Range("I27").Value = 2229.4167
$ws.Range("J27").Value = 3550.7856
$ws.Range("K27").Value = 2229.4167
$ws.Range("L27").Value = 3550.7856
$ws.Range("M27").Value = -2122.4167
$ws.Range("N27").Value = -3764.7856

$ws.Range("H68").Value = 2842525.8
$ws.Range("I68").Value = 5683018
$ws.Range("K68").Value = 5683018
$ws.Range("M68").Value = -5682269

$ws.Range("H71").Value = 2842525.8
$ws.Range("I71").Value = 5683018
$ws.Range("K71").Value = 28415090
$ws.Range("M71").Value = -28411346

$ws.Range("H136").Value = 3952.3572
$ws.Range("I136").Value = 3972.198
$ws.Range("K136").Value = 11916.594
$ws.Range("M136").Value = -9366.593999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 15136.333
$ws.Range("J45").Value = 13204.5
$ws.Range("L45").Value = 13204.5
$ws.Range("N45").Value = -14186.5

$ws.Range("H132").Value = 17862394
$ws.Range("I132").Value = 6322.6875
$ws.Range("J132").Value = 41670490
$ws.Range("K132").Value = 18968.0625
$ws.Range("L132").Value = 125011470
$ws.Range("M132").Value = -16438.0625
$ws.Range("N132").Value = -125016530

$ws.Range("H136").Value = 5621.9478
$ws.Range("I136").Value = 2663.1086
$ws.Range("J136").Value = 8344.08
$ws.Range("K136").Value = 7989.325800000001
$ws.Range("L136").Value = 25032.24
$ws.Range("M136").Value = -5439.325800000001
$ws.Range("N136").Value = -30132.24
